# Applies the weekly Fruta/hortaliza refresh for the
# "Fruta, Vega Modelo de Temuco - Maracuya" sheet: existing rows 2-24
# are reshuffled/updated in place (Fecha, Volumen, Precio min/max/prom,
# Unidad de comercializacion, Origen, Precio $/Kg, Kg/unidad).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44448
$ws.Range("M2").Value = 50
$ws.Range("N2").Value = 38000
$ws.Range("O2").Value = 38000
$ws.Range("P2").Value = 38000
$ws.Range("R2").Value = 'Región de Arica y Parinacota'
$ws.Range("S2").Value = 2111
# Row 3
$ws.Range("D3").Value = 44364
$ws.Range("M3").Value = 90
$ws.Range("N3").Value = 1700
$ws.Range("O3").Value = 1700
$ws.Range("P3").Value = 1700
$ws.Range("Q3").Value = '$/kilo'
$ws.Range("S3").Value = 1700
$ws.Range("T3").Value = 1
# Row 4
$ws.Range("D4").Value = 44379
$ws.Range("M4").Value = 10
$ws.Range("N4").Value = 30000
$ws.Range("O4").Value = 30000
$ws.Range("P4").Value = 30000
$ws.Range("R4").Value = 'Región de Arica y Parinacota'
$ws.Range("S4").Value = 1667
# Row 5
$ws.Range("D5").Value = 44432
$ws.Range("M5").Value = 10
$ws.Range("R5").Value = 'Perú'
# Row 6
$ws.Range("D6").Value = 44434
$ws.Range("M6").Value = 40
$ws.Range("N6").Value = 35000
$ws.Range("O6").Value = 35000
$ws.Range("P6").Value = 35000
$ws.Range("S6").Value = 1944
# Row 7
$ws.Range("D7").Value = 44449
$ws.Range("M7").Value = 20
$ws.Range("R7").Value = 'Región de Arica y Parinacota'
# Row 8
$ws.Range("D8").Value = 44405
$ws.Range("M8").Value = 10
$ws.Range("N8").Value = 35000
$ws.Range("O8").Value = 35000
$ws.Range("P8").Value = 35000
$ws.Range("Q8").Value = '$/caja 18 kilos'
$ws.Range("S8").Value = 1944
$ws.Range("T8").Value = 18
# Row 9
$ws.Range("D9").Value = 44377
$ws.Range("M9").Value = 30
$ws.Range("N9").Value = 40000
$ws.Range("O9").Value = 40000
$ws.Range("P9").Value = 40000
$ws.Range("S9").Value = 2222
# Row 10
$ws.Range("D10").Value = 44424
$ws.Range("M10").Value = 15
$ws.Range("N10").Value = 35000
$ws.Range("O10").Value = 35000
$ws.Range("P10").Value = 35000
$ws.Range("S10").Value = 1944
# Row 11
$ws.Range("D11").Value = 44369
$ws.Range("M11").Value = 5
$ws.Range("R11").Value = 'Perú'
# Row 14
$ws.Range("D14").Value = 44442
$ws.Range("M14").Value = 15
$ws.Range("R14").Value = 'Perú'
# Row 15
$ws.Range("D15").Value = 44392
$ws.Range("M15").Value = 20
# Row 17
$ws.Range("D17").Value = 44363
$ws.Range("M17").Value = 144
# Row 18
$ws.Range("D18").Value = 44438
$ws.Range("M18").Value = 25
$ws.Range("O18").Value = 35000
$ws.Range("P18").Value = 35000
$ws.Range("S18").Value = 1944
# Row 20
$ws.Range("D20").Value = 44357
$ws.Range("M20").Value = 10
$ws.Range("N20").Value = 38000
$ws.Range("O20").Value = 38000
$ws.Range("P20").Value = 38000
$ws.Range("S20").Value = 2111
# Row 21
$ws.Range("D21").Value = 44294
$ws.Range("M21").Value = 15
$ws.Range("N21").Value = 35000
$ws.Range("O21").Value = 35000
$ws.Range("P21").Value = 35000
$ws.Range("S21").Value = 1944
# Row 22
$ws.Range("D22").Value = 44264
$ws.Range("M22").Value = 20
$ws.Range("N22").Value = 40000
$ws.Range("O22").Value = 40000
$ws.Range("P22").Value = 40000
$ws.Range("S22").Value = 2222
# Row 23
$ws.Range("D23").Value = 44279
$ws.Range("M23").Value = 30
$ws.Range("N23").Value = 35000
$ws.Range("O23").Value = 36000
$ws.Range("P23").Value = 35667
$ws.Range("S23").Value = 1982
# Row 24
$ws.Range("D24").Value = 44418
$ws.Range("M24").Value = 30
